# "planning changed and costs added"
# Restructure the purchase-planning sheet: move the "what/link/Price/quantity"
# table from columns A/C/L/M into a compact A:E table, add two new purchase
# rows (Micro USB OTG, LiPo Batteries + Charger), add per-row subtotal
# formulas and a grand Total row, and tidy up number formats/column widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Grab the text of cells that must be preserved verbatim (incl. the
#     hyperlink-display string that still carries a trailing space) before we
#     clear the sheet, so the shared-string table keeps reusing the same
#     entries instead of minting near-duplicates. ---
$txtWhat = $ws.Range("A1").Value2
$txtLink = $ws.Range("C1").Value2
$txtPrice = $ws.Range("L1").Value2
$txtQty = $ws.Range("M1").Value2
$txtDrone = $ws.Range("A2").Value2
$txtDroneUrl = $ws.Range("C2").Value2

# --- Wipe the old layout (cells + hyperlinks) so we can rebuild cleanly. ---
$ws.Cells.Clear()
$ws.Hyperlinks.Delete()

# --- Column widths for the new A/B columns. ---
$ws.Columns("A").ColumnWidth = 16.5
$ws.Columns("B").ColumnWidth = 22

# --- Header row (bold). ---
$ws.Range("A1").Value = $txtWhat
$ws.Range("B1").Value = $txtLink
$ws.Range("C1").Value = $txtPrice
$ws.Range("D1").Value = $txtQty
$ws.Range("E1").Value = "Total"
$ws.Range("A1:E1").Font.Bold = $true

# --- Row 2: Parrot AR.drone 2.0 (kept) ---
$ws.Range("A2").Value = $txtDrone
$ws.Range("B2").Value = $txtDroneUrl
$ws.Range("C2").Value = 116.71
$ws.Range("D2").Value = 1

# --- Row 3: Micro USB OTG (new; styled like a link but no real hyperlink,
#     matching the source workbook) ---
$ws.Range("A3").Value = "Micro USB OTG "
$ws.Range("B3").Value = "https://nl.aliexpress.com/item/1pc-Micro-USB-Type-B-Male-To-Micro-B-Male-5-Pin-Converter-OTG-Adapter-Lead/32808154910.html?spm=a2g0z.search0104.3.221.23ff22089TknZw&ws_ab_test=searchweb0_0,searchweb201602_1_10152_10151_10065_10344_10068_10342_10343_10340_10341_10084_10083_10618_10630_10304_10307_10302_5711211_10313_10059_5722311_10534_100031_10103_10627_10626_10624_10623_10622_10621_10620_5711311_10142,searchweb201603_25,ppcSwitch_4&algo_expid=24fc4199-da5a-4461-8ea1-e398eb22a66b-29&algo_pvid=24fc4199-da5a-4461-8ea1-e398eb22a66b&priceBeautifyAB=0"
$ws.Range("B3").Style = "Hyperlink"
$ws.Range("C3").Value = 1.32
$ws.Range("D3").Value = 2

# --- Row 4: LiPo Batteries + Charger (new, with real hyperlink) ---
$ws.Range("A4").Value = "LiPo Batteries + Charger"
$ws.Range("B4").Value = "https://nl.aliexpress.com/item/5pcs-Lipo-Battery-3-7V-150mAh-USB-Lipo-Charger-Set-For-JJRC-H20-Mini-RC-Hexacopter/32791721090.html?spm=a2g0s.13010208.99999999.262.0rynF9"
$ws.Range("C4").Value = 12.47
$ws.Range("D4").Value = 1

# --- Quantity column: integer number format ---
$ws.Range("D2:D4").NumberFormat = "0"

# --- Price / subtotal currency format (euro, no red-negative variant) ---
$ws.Range("C2:C4").NumberFormat = '"€"\ #,##0.00'

# --- Per-row subtotal formulas (Price * quantity); E3:E4 become one shared
#     formula block, matching how Excel would fill this down. ---
$ws.Range("E2").Formula = "=C2*D2"
$ws.Range("E3:E4").Formula = "=C3*D3"
$ws.Range("E2:E4").NumberFormat = '"€"\ #,##0.00'

# --- Grand total row ---
$ws.Range("D5").Value = "Total"
$ws.Range("D5").Font.Bold = $true
$ws.Range("E5").Formula = "=SUM(E2:E4)"
$ws.Range("E5").NumberFormat = '"€"\ #,##0.00'
$ws.Range("E5").Font.Bold = $true

# --- Real hyperlinks (Amazon link reused on B2, new AliExpress link on B4).
#     Adding the hyperlink with a display string overwrites the cell's text,
#     so we pass the trimmed text to get the right `display=` attribute and
#     then restore the original (trailing-space) text into the cell
#     afterwards, which also makes it reuse the existing shared-string entry. ---
$ws.Hyperlinks.Add($ws.Range("B2"), $txtDroneUrl.Trim(), "", "", $txtDroneUrl.Trim())
$ws.Range("B2").Value = $txtDroneUrl
$ws.Range("B2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("B4"), "https://nl.aliexpress.com/item/5pcs-Lipo-Battery-3-7V-150mAh-USB-Lipo-Charger-Set-For-JJRC-H20-Mini-RC-Hexacopter/32791721090.html?spm=a2g0s.13010208.99999999.262.0rynF9")
$ws.Range("B4").Style = "Hyperlink"

# --- Selection left where the author left it ---
$ws.Range("D6").Select()
